# Add data for 2022-03-25: one more day's worth of carjackings rolled into
# the "current month" column (column B), moving the "through" date from
# March 16 to March 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet tab name and the column-B header label to reflect the
# new "through" date.
$ws.Name = "Through 2022-03-17"
$ws.Range("B1").Value = "March 2022 (through March 17)"

# Incremented existing counts.
$ws.Range("E4").Value = 8
$ws.Range("N4").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("W5").Value = 4
$ws.Range("B9").Value = 4
$ws.Range("H44").Value = 3
$ws.Range("N49").Value = 2
$ws.Range("K77").Value = 2

# Brand-new counts (previously-empty cells).
$ws.Range("H3").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("H67").Value = 1
$ws.Range("B68").Value = 1
$ws.Range("B74").Value = 1
